# Fix mis-encoded "±" (plus-minus) characters that were double-encoded as
# "Â±" (UTF-8 bytes interpreted as Latin-1) back into the correct "±" glyph.
# Affects columns B (score), C (training_time), D (test_time) for the data
# rows 2-17 on the active worksheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The original text was UTF-8 encoded bytes for "±" (0xC2 0xB1) that got
# re-interpreted as Latin-1/Windows-1252 and re-encoded, producing the two
# characters U+00C2 (Â) followed by U+00B1 (±) in place of the single
# intended U+00B1 (±) character.
$badChar  = [string][char]0xC2 + [string][char]0xB1
$goodChar = [string][char]0xB1

for ($row = 2; $row -le 17; $row++) {
    foreach ($col in @("B", "C", "D")) {
        $cell = $ws.Range("$col$row")
        $value = $cell.Value2
        if ($value -ne $null -and $value.Contains($badChar)) {
            $cell.Value = $value.Replace($badChar, $goodChar)
        }
    }
}
